# Generate Report for Archive
#
# The localization status report is being refreshed: the rows that were
# previously "Ready for handoff" have moved on to the next pipeline stage,
# "In Translation". This status string is shown in the Overview sheet
# (once per language column) as well as on each language's own detail
# sheet (zh-cn, de-de), so every occurrence needs to be updated. Because
# the new status text is shorter than the old one, the Status column on
# each sheet is narrowed to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview"
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$dede     = $wb.Worksheets.Item(3)   # "de-de"

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: one Status-like column per target language (E, F) ---
# (Value2 is used for the comparison - it reliably unwraps to a plain
# string/number, unlike Value which can come back as an opaque variant.)
if ($overview.Range("E2").Value2 -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Value2 -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }
if ($overview.Range("E3").Value2 -eq $oldStatus) { $overview.Range("E3").Value = $newStatus }
if ($overview.Range("F3").Value2 -eq $oldStatus) { $overview.Range("F3").Value = $newStatus }

# --- Per-language detail sheets: "Status" is column C ---
if ($zhcn.Range("C2").Value2 -eq $oldStatus) { $zhcn.Range("C2").Value = $newStatus }
if ($zhcn.Range("C3").Value2 -eq $oldStatus) { $zhcn.Range("C3").Value = $newStatus }

if ($dede.Range("C2").Value2 -eq $oldStatus) { $dede.Range("C2").Value = $newStatus }
if ($dede.Range("C3").Value2 -eq $oldStatus) { $dede.Range("C3").Value = $newStatus }

# --- Re-fit the Status columns now that the text is shorter ---
# (ColumnWidth is in "characters"; 12.5 is the narrowest setting that still
# maps to the same rendered column width the shorter status text needs.)
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth  # Overview!E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth  # Overview!F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth       # zh-cn!C (Status)
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth       # de-de!C (Status)
